$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (Price / Volume(1h)) per refreshed coinranking.com snapshot.
# Force Text number format before assigning so values like "278.00" and "0.99%" are
# stored as literal strings (matching the existing inlineStr text cells) rather than
# being auto-converted to numbers/percentages by Excel's input parser.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.99%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.51%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.874"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.01%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06432"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.48%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.007"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.42%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.192"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.48%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8842"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.18%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1560"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.75%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05109"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.67%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07490"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.10%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02886"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.41%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08973"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.80%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001575"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.25%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006401"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.59%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006129"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.86%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.94%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.309"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.18%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.07%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.93%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.913"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.08%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04423"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.04%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001174"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.49%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003876"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.01%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.78%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-1.63%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04147"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006808"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-2.45%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.37%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001889"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-11.86%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01129"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.63%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005324"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.70%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "13.27%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01853"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-11.72%"
